$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-ListParagraphXml([string]$text) {
    return '<w:p ' + $wNs + '>' +
             '<w:pPr>' +
               '<w:pStyle w:val="Listenabsatz"/>' +
               '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
               '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
             '</w:pPr>' +
             '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + $text + '</w:t></w:r>' +
           '</w:p>'
}

# --- Locate the three paragraphs that need to be replaced -----------------
# 1) the empty paragraph right before "Find key in yard"
# 2) the "Find key in yard" list paragraph
# 3) the "Go into hospital" list paragraph (holds the _GoBack bookmark)
$findKeyIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Find key in yard*") {
        $findKeyIdx = $i
        break
    }
}
if ($findKeyIdx -lt 0) {
    throw "Could not locate the 'Find key in yard' paragraph"
}

$emptyPara = $d.Paragraphs.Item($findKeyIdx - 1)

# --- Step 1: turn the lone empty paragraph into the first new list item ---
# (Doing this as its own single-paragraph replace avoids an engine quirk
#  where InsertXML on a multi-paragraph range leaves a *leading* empty
#  paragraph untouched instead of replacing it.)
[void]$emptyPara.Range.InsertXML((New-ListParagraphXml("Light goes out")))

# --- Step 2: replace "Find key in yard" / "Go into hospital" --------------
# Paragraph indices are unchanged since step 1 was a strict 1-for-1 swap.
$findKeyPara     = $d.Paragraphs.Item($findKeyIdx)
$goHospitalPara  = $d.Paragraphs.Item($findKeyIdx + 1)

$xml = ""
$xml += New-ListParagraphXml("Find lamp")
$xml += New-ListParagraphXml("Find key")
$xml += New-ListParagraphXml("Go inside")
$xml += New-ListParagraphXml("Turn on power source")

# final (previously "Go into hospital") paragraph: keeps the style but loses
# the numbering and text, and still carries the _GoBack bookmark.
$xml += '<w:p ' + $wNs + '>' +
             '<w:pPr>' +
               '<w:pStyle w:val="Listenabsatz"/>' +
               '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
             '</w:pPr>' +
             '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
             '<w:bookmarkEnd w:id="0"/>' +
           '</w:p>'

$target = $d.Range($findKeyPara.Range.Start, $goHospitalPara.Range.End)
[void]$target.InsertXML($xml)
